$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New D (Fecha), M (Volumen), N (Precio minimo), O (Precio maximo),
# P (Precio promedio ponderado), S (Precio $/Kg) values for rows 2-22.
$data = @(
    @{Row=2;  D=44431; M=100; N=1300; O=1300; P=1300; S=1300},
    @{Row=3;  D=44748; M=300; N=2300; O=2300; P=2300; S=2300},
    @{Row=4;  D=44473; M=120; N=1200; O=1200; P=1200; S=1200},
    @{Row=5;  D=44435; M=130; N=1300; O=1300; P=1300; S=1300},
    @{Row=6;  D=44830; M=50;  N=2500; O=2500; P=2500; S=2500},
    @{Row=7;  D=44424; M=50;  N=1200; O=1200; P=1200; S=1200},
    @{Row=8;  D=44760; M=80;  N=2300; O=2300; P=2300; S=2300},
    @{Row=9;  D=44417; M=80;  N=1200; O=1200; P=1200; S=1200},
    @{Row=10; D=44405; M=50;  N=1200; O=1200; P=1200; S=1200},
    @{Row=11; D=44811; M=60;  N=2500; O=2500; P=2500; S=2500},
    @{Row=12; D=44343; M=60;  N=1300; O=1300; P=1300; S=1300},
    @{Row=13; D=44749; M=120; N=2300; O=2300; P=2300; S=2300},
    @{Row=14; D=44432; M=30;  N=1300; O=1300; P=1300; S=1300},
    @{Row=15; D=44753; M=160; N=2300; O=2300; P=2300; S=2300},
    @{Row=16; D=44812; M=50;  N=2500; O=2500; P=2500; S=2500},
    @{Row=17; D=44476; M=80;  N=1200; O=1200; P=1200; S=1200},
    @{Row=18; D=44418; M=40;  N=1200; O=1200; P=1200; S=1200},
    @{Row=19; D=44438; M=60;  N=1200; O=1200; P=1200; S=1200},
    @{Row=20; D=44762; M=50;  N=2300; O=2300; P=2300; S=2300},
    @{Row=21; D=44763; M=50;  N=2300; O=2300; P=2300; S=2300},
    @{Row=22; D=44357; M=35;  N=1000; O=1000; P=1000; S=1000}
)

foreach ($entry in $data) {
    $r = $entry.Row
    $ws.Cells.Item($r, 4).Value = $entry.D    # D
    $ws.Cells.Item($r, 13).Value = $entry.M   # M
    $ws.Cells.Item($r, 14).Value = $entry.N   # N
    $ws.Cells.Item($r, 15).Value = $entry.O   # O
    $ws.Cells.Item($r, 16).Value = $entry.P   # P
    $ws.Cells.Item($r, 19).Value = $entry.S   # S
}
